# Applies the "a few tweaks and renaming" commit to data.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) battery_level (column K) bumped from 70 to 90 for every row that had 70,
#    i.e. every data row except rows 23-28 (which hold 8 / 40 and are left
#    untouched).
# ---------------------------------------------------------------------------
$kRowsToBump = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43)
foreach ($r in $kRowsToBump) {
    $ws.Cells.Item($r, 11).Value = 90
}

# ---------------------------------------------------------------------------
# 2) follower_wellbeing (column O) tweaks
# ---------------------------------------------------------------------------
$ws.Range("O18").Value = 0.26
$ws.Range("O19").Value = 0.26
$ws.Range("O21").Value = -0.03
$ws.Range("O22").Value = -0.03
$ws.Range("O27").Value = -0.03
$ws.Range("O28").Value = -0.03

# ---------------------------------------------------------------------------
# 3) not_follow_locations (column J) - row 29 had the stray value "s",
#    it should read "[]" like the other non-follow cases.
# ---------------------------------------------------------------------------
$ws.Range("J29").Value = "[]"

# ---------------------------------------------------------------------------
# 4) row 38 tweaks
# ---------------------------------------------------------------------------
$ws.Range("D38").Value = 3
$ws.Range("R38").Value = 0

# ---------------------------------------------------------------------------
# 5) rows 39-43 - the follow-up "go_to_last_seen / stay / go_to_charge_station"
#    cases got reshuffled and a new row 44 was appended.
# ---------------------------------------------------------------------------
$ws.Range("N39").Value = -0.7
$ws.Range("Q39").Value = "go_to_last_seen"
$ws.Range("R39").Value = 1

$ws.Range("Q40").Value = "stay"

$ws.Range("D41").Value = 10
$ws.Range("E41").Value = 1
$ws.Range("F41").Value = 2
$ws.Range("N41").Value = 1
$ws.Range("Q41").Value = "go_to_charge_station"
$ws.Range("R41").Value = 0

$ws.Range("N42").Value = -0.7
$ws.Range("Q42").Value = "go_to_last_seen"
$ws.Range("R42").Value = 1

$ws.Range("Q43").Value = "stay"

# New row 44 (case_id 43) - holds what used to be in row 41 before the
# reshuffle above (20 / 0.5 / 0, go_to_charge_station).
$ws.Range("A44").Value = 43
$ws.Range("B44").Formula = "=FALSE()"
$ws.Range("D44").Value = 20
$ws.Range("E44").Value = 0.5
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = "bathroom"
$ws.Range("H44").Value = "bedroom"
$ws.Range("I44").Formula = "=TRUE()"
$ws.Range("J44").Value = "['bathroom']"
$ws.Range("K44").Value = 90
$ws.Range("L44").Value = "[]"
$ws.Range("M44").Value = "day"
$ws.Range("N44").Value = 1
$ws.Range("O44").Value = 0.93
$ws.Range("P44").Value = 0.98
$ws.Range("Q44").Value = "go_to_charge_station"
$ws.Range("R44").Value = 0
$ws.Range("S44").Value = "wellbeing"
$ws.Range("B44").NumberFormat = "General"
$ws.Range("I44").NumberFormat = "General"

# ---------------------------------------------------------------------------
# 6) fix mangled apostrophes: "[‘move_away’]" -> "['move_away']"
# ---------------------------------------------------------------------------
$fixedInstr = "['move_away']"
$moveAwayRows = @(11,12,13,29,30,31,32,33,34,35,36,37)
foreach ($r in $moveAwayRows) {
    $ws.Cells.Item($r, 12).Value = $fixedInstr
}

# ---------------------------------------------------------------------------
# 7) drop the now-redundant cell style (numFmtId 164, duplicate of style 0)
#    that used to be applied to every not_follow_request / seen boolean cell
#    in rows 3-43; normalising the number format collapses them back onto
#    the shared "General" style.
# ---------------------------------------------------------------------------
$ws.Range("I3:I43").NumberFormat = "General"
$ws.Range("B15:B43").NumberFormat = "General"

# ---------------------------------------------------------------------------
# 8) view state - the workbook was last left scrolled to C4 with D41 selected
# ---------------------------------------------------------------------------
$ws.Range("D41").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 4
